$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, shifting rows 72:148 down to 73:149.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new data.
$ws.Cells.Item(72, 1).Value = 11
$ws.Cells.Item(72, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(72, 3).Value = "Bíobío"
$ws.Cells.Item(72, 4).Value = 44539
$ws.Cells.Item(72, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(72, 5).Value = 8
$ws.Cells.Item(72, 6).Value = 100114001
$ws.Cells.Item(72, 7).Value = "Papa"
$ws.Cells.Item(72, 8).Value = "Asterix"
$ws.Cells.Item(72, 9).Value = "1a (cosecha)"
$ws.Cells.Item(72, 10).Value = 270
$ws.Cells.Item(72, 11).Value = 9500
$ws.Cells.Item(72, 12).Value = 10000
$ws.Cells.Item(72, 13).Value = 9722
$ws.Cells.Item(72, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(72, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(72, 16).Value = 389
$ws.Cells.Item(72, 17).Value = 25
$ws.Cells.Item(72, 18).Value = "Hortaliza"
